$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.520.71'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '1.651.85'
$ws.Range('E3').Value = '  -3.36%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9993'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3648'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.57%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '46.36'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3249'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.124'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.63%  '
$ws.Range('E11').Value = '  -6.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9991'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.956'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.36'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -8.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.598'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.01%  '
$ws.Range('D16').Value = '1.655.56'
$ws.Range('E16').Value = '  -2.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001041'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06583'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.96%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '78.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.925'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.53'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.53%  '
$ws.Range('D24').Value = '24.498.03'
$ws.Range('E24').Value = '  -1.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.463'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.330'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -16.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '146.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -9.33%  '
$ws.Range('D29').Value = '1.833.87'
$ws.Range('E29').Value = '  -3.22%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.187'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.58%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '124.03'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.068'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.705'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -16.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08435'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.15%  '
$ws.Range('E35').Value = '  -6.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -12.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.269'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.191'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06009'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02220'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.14%  '
$ws.Range('E41').Value = '  -7.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.080'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -11.93%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5884'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.801'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.67'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5608'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '123.71'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.942'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06921'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.188'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.91%  '
